$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# dS0 (column E) correction
$ws.Range("E9").Value = 3

# dSF (column F) repull / recalculation of values
$ws.Range("F2").Value = 3
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = -1
$ws.Range("F10").Value = 0
$ws.Range("F12").Value = 2
$ws.Range("F13").Value = 4
$ws.Range("F19").Value = 2
$ws.Range("F20").Value = 1
$ws.Range("F21").Value = 2
$ws.Range("F27").Value = -7
$ws.Range("F28").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("F38").Value = 2
$ws.Range("F39").Value = 1
$ws.Range("F42").Value = -4
$ws.Range("F47").Value = 0
$ws.Range("F61").Value = 2
$ws.Range("F64").Value = 0
$ws.Range("F67").Value = 1
